$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto values scraped on Wed Nov  8 05:46:16 UTC 2023

# Row 2
$ws.Range("D2").Value = "'35.303.40"
$ws.Range("E2").Value = '  +0.54%  '

# Row 3
$ws.Range("D3").Value = "'1.880.09"
$ws.Range("E3").Value = '  -1.30%  '

# Row 4
$ws.Range("E4").Value = '  -0.60%  '

# Row 5
$ws.Range("D5").Value = "'245.13"
$ws.Range("E5").Value = '  -3.26%  '

# Row 6
$ws.Range("D6").Value = "'0.684"
$ws.Range("E6").Value = '  -1.61%  '

# Row 7
$ws.Range("E7").Value = '  -0.64%  '

# Row 8
$ws.Range("D8").Value = "'43.33"
$ws.Range("E8").Value = '  +4.68%  '

# Row 9
$ws.Range("E9").Value = '  -1.43%  '

# Row 10
$ws.Range("D10").Value = "'53.37"
$ws.Range("E10").Value = '  +1.63%  '

# Row 11
$ws.Range("D11").Value = "'0.0739"
$ws.Range("E11").Value = '  -1.59%  '

# Row 12
$ws.Range("D12").Value = "'0.0972"
$ws.Range("E12").Value = '  -0.92%  '

# Row 13
$ws.Range("D13").Value = "'13.47"
$ws.Range("E13").Value = '  +1.96%  '

# Row 14
$ws.Range("D14").Value = "'2.153.67"
$ws.Range("E14").Value = '  -1.33%  '

# Row 15
$ws.Range("D15").Value = "'0.769"
$ws.Range("E15").Value = '  +4.82%  '

# Row 16
$ws.Range("B16").Value = "'Polkadot"
$ws.Range("C16").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'4.90"
$ws.Range("E16").Value = '  -2.24%  '

# Row 17
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'1.871.81"
$ws.Range("E17").Value = '  -1.88%  '

# Row 18
$ws.Range("D18").Value = "'35.299.38"
$ws.Range("E18").Value = '  +0.51%  '

# Row 19
$ws.Range("D19").Value = "'72.78"
$ws.Range("E19").Value = '  -1.29%  '

# Row 20
$ws.Range("E20").Value = '  -2.23%  '

# Row 21
$ws.Range("D21").Value = "'243.62"
$ws.Range("E21").Value = '  +0.26%  '

# Row 22
$ws.Range("E22").Value = '  -1.50%  '

# Row 23
$ws.Range("E23").Value = '  -2.00%  '

# Row 24
$ws.Range("D24").Value = "'2.63"
$ws.Range("E24").Value = '  +8.15%  '

# Row 25
$ws.Range("E25").Value = '  -0.62%  '

# Row 26
$ws.Range("D26").Value = "'2.16"
$ws.Range("E26").Value = '  -6.27%  '

# Row 27
$ws.Range("D27").Value = "'165.35"
$ws.Range("E27").Value = '  -1.56%  '

# Row 28
$ws.Range("D28").Value = "'8.54"
$ws.Range("E28").Value = '  -0.54%  '

# Row 29
$ws.Range("D29").Value = "'18.26"
$ws.Range("E29").Value = '  -1.55%  '

# Row 30
$ws.Range("E30").Value = '  -2.33%  '

# Row 31
$ws.Range("D31").Value = "'4.128.46"
$ws.Range("E31").Value = '  -0.01%  '

# Row 32
$ws.Range("E32").Value = '  +8.01%  '

# Row 33
$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.26"
$ws.Range("E33").Value = '  -1.98%  '

# Row 34
$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0586"
$ws.Range("E34").Value = '  -2.94%  '

# Row 35
$ws.Range("B35").Value = "'WEMIXToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.96"
$ws.Range("E35").Value = '  -2.85%  '

# Row 36
$ws.Range("D36").Value = "'4.13"
$ws.Range("E36").Value = '  -2.13%  '

# Row 37
$ws.Range("E37").Value = '  -0.61%  '

# Row 38
$ws.Range("E38").Value = '  -1.53%  '

# Row 39
$ws.Range("B39").Value = "'Kaspa"
$ws.Range("C39").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.0732"
$ws.Range("E39").Value = '  +12.37%  '

# Row 40
$ws.Range("B40").Value = "'LidoDAOToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'1.93"
$ws.Range("E40").Value = '  -3.85%  '

# Row 41
$ws.Range("E41").Value = '  +0.70%  '

# Row 42
$ws.Range("E42").Value = '  +0.19%  '

# Row 43
$ws.Range("D43").Value = "'96.10"
$ws.Range("E43").Value = '  -7.32%  '

# Row 44
$ws.Range("E44").Value = '  -2.67%  '

# Row 45
$ws.Range("D45").Value = "'1.303.80"
$ws.Range("E45").Value = '  -0.25%  '

# Row 46
$ws.Range("E46").Value = '  -2.70%  '

# Row 47
$ws.Range("D47").Value = "'0.0797"
$ws.Range("E47").Value = '  +6.58%  '

# Row 48
$ws.Range("E48").Value = '  -1.75%  '

# Row 49
$ws.Range("E49").Value = '  -1.25%  '

# Row 50
$ws.Range("D50").Value = "'12.13"
$ws.Range("E50").Value = '  -4.19%  '

# Row 51
$ws.Range("D51").Value = "'6.22"
$ws.Range("E51").Value = '  -5.64%  '
